# Updated symbol list on Wed Jan 25 23:29:41 UTC 2023 with GitHub Actions
# Refreshes price/volume readings for the coinranking.com scrape and
# re-ranks the 7-17 block (GateToken moves to rank #1, pushing the rest
# down by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '307.32'
    'E2' = '1.27%'
    'D3' = '36.53'
    'E3' = '3.62%'
    'D4' = '5.092'
    'E4' = '1.64%'
    'D5' = '0.08095'
    'E5' = '3.26%'
    'D6' = '1.974'
    'E6' = '4.43%'
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D7' = '4.126'
    'E7' = '0.66%'
    'B8' = 'KuCoinToken'
    'C8' = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
    'D8' = '7.763'
    'E8' = '-0.68%'
    'B9' = 'MXToken'
    'C9' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D9' = '0.9354'
    'E9' = '1.56%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.1452'
    'E10' = '35.77%'
    'B11' = 'WazirX'
    'C11' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D11' = '0.1931'
    'E11' = '3.11%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.09144'
    'E12' = '-1.87%'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '0.03542'
    'E13' = '-3.83%'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '0.09792'
    'E14' = '-1.09%'
    'B15' = 'BitForexToken'
    'C15' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D15' = '0.001432'
    'E15' = '1.35%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.005809'
    'E16' = '1.70%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.521'
    'E17' = '1.74%'
    'E18' = '4.82%'
    'D19' = '0.3432'
    'E19' = '-0.01%'
    'D20' = '0.1302'
    'E20' = '0.21%'
    'D21' = '4.955'
    'E21' = '-3.18%'
    'D22' = '0.2402'
    'E22' = '9.17%'
    'D23' = '0.04505'
    'E23' = '-0.99%'
    'D24' = '0.001212'
    'E24' = '-1.39%'
    'D25' = '0.004859'
    'E25' = '4.18%'
    'D26' = '0.0001241'
    'E26' = '-1.16%'
    'D27' = '0.0004444'
    'E27' = '-0.50%'
    'D39' = '0.01976'
    'E39' = '4.95%'
    'D40' = '0.04872'
    'E40' = '3.19%'
    'D41' = '0.01103'
    'E41' = '10.23%'
    'D42' = '0.007556'
    'E42' = '-0.45%'
    'D43' = '0.1367'
    'E43' = '2.51%'
    'D44' = '0.002111'
    'E44' = '-0.58%'
    'D45' = '0.009766'
    'E45' = '-13.52%'
    'D46' = '0.00006375'
    'E46' = '1.37%'
    'D47' = '0.00000000750'
    'E47' = '-0.37%'
    'E48' = '0.43%'
    'D49' = '0.001191'
    'E49' = '-8.76%'
    'D50' = '0.00002101'
    'E50' = '-0.37%'
    'D51' = '0.0002001'
    'E51' = '-0.37%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (prices, "x.xx%"
    # volumes) keep their exact original formatting (trailing zeros,
    # "%" suffix, long decimal strings) instead of being coerced into
    # floats/percentages by COM's type inference.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop back to the default style so we don't leave a stray
    # text-format style attached to the cell (matches the source file,
    # where these cells carry no explicit style index).
    $cell.Style = "Normal"
}
